$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 6, shifting the old row 6 down to row 7
$ws.Rows.Item(6).Insert()

# Fill in the new row 6 with the data
$ws.Range("A6").Value = "djfkad"
$ws.Range("B6").Value = 98283

# Update the selection to match the new active cell/range
$ws.Activate()
$ws.Range("A6:B6").Select()
